$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178..235 down to 179..236.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new record.
$ws.Range("A178").Value = 8
$ws.Range("B178").Value = "Terminal La Palmera de La Serena"
$ws.Range("C178").Value = "Coquimbo"
$ws.Range("D178").Value = 44559
$ws.Range("E178").Value = 4
$ws.Range("F178").Value = 100112032
$ws.Range("G178").Value = "Zapallo italiano"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 600
$ws.Range("K178").Value = 7000
$ws.Range("L178").Value = 8000
$ws.Range("M178").Value = 7500
$ws.Range("N178").Value = '$/caja 70 unidades'
$ws.Range("O178").Value = 'Provincia de Limarí'
$ws.Range("P178").Value = 107
$ws.Range("Q178").Value = 70
$ws.Range("R178").Value = "Hortaliza"
